$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting existing rows 79..141 down to 80..142
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new record
$ws.Cells.Item(79, 1).Value = 11
$ws.Cells.Item(79, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(79, 3).Value = "Bíobío"
$ws.Cells.Item(79, 4).Value = 44873
$ws.Cells.Item(79, 5).Value = 8
$ws.Cells.Item(79, 6).Value = 100112021
$ws.Cells.Item(79, 7).Value = "Ají"
$ws.Cells.Item(79, 8).Value = "Inferno"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 80
$ws.Cells.Item(79, 11).Value = 16000
$ws.Cells.Item(79, 12).Value = 17000
$ws.Cells.Item(79, 13).Value = 16375
$ws.Cells.Item(79, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 1638
$ws.Cells.Item(79, 17).Value = 10
$ws.Cells.Item(79, 18).Value = "Hortaliza"
